# Update "想去人数" (number of people interested) counts across sheets.
# Source data sheets: 展览 (Exhibition) and 演出 (Show).
# 全部类型 (All types) is a combined/aggregate sheet that mirrors the same
# rows, so it must receive the matching updates as well.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 23
$ws1.Range("F7").Value  = 1725
$ws1.Range("F11").Value = 1815
$ws1.Range("F18").Value = 17
$ws1.Range("F22").Value = 775
$ws1.Range("F26").Value = 262

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 7

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 23
$ws4.Range("F7").Value  = 1725
$ws4.Range("F8").Value  = 7
$ws4.Range("F12").Value = 1815
$ws4.Range("F19").Value = 17
$ws4.Range("F23").Value = 775
$ws4.Range("F27").Value = 262
